$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete now-obsolete trailing rows 54:56 (table shrinks from 56 to 53 data-ish rows)
$ws.Range("A54:B56").EntireRow.Delete()

# Rewrite sector labels + correlation values for rows 2-53 (sector order/counts changed)
$ws.Cells.Item(2, 1).Value = "Industrial Conglomerates(6)"
$ws.Cells.Item(2, 2).Value = 0.7121166213872037
$ws.Cells.Item(3, 1).Value = "Road & Rail(22)"
$ws.Cells.Item(3, 2).Value = 0.5989653154570562
$ws.Cells.Item(4, 1).Value = "Air Freight & Logistics(11)"
$ws.Cells.Item(4, 2).Value = 0.5775675828705409
$ws.Cells.Item(5, 1).Value = "Energy Equipment & Services(32)"
$ws.Cells.Item(5, 2).Value = 0.5751299603686909
$ws.Cells.Item(6, 1).Value = "Marine(15)"
$ws.Cells.Item(6, 2).Value = 0.5690460143460049
$ws.Cells.Item(7, 1).Value = "Construction & Engineering(20)"
$ws.Cells.Item(7, 2).Value = 0.5675362280049581
$ws.Cells.Item(8, 1).Value = "Construction Materials(8)"
$ws.Cells.Item(8, 2).Value = 0.5663544204918297
$ws.Cells.Item(9, 1).Value = "Trading Companies & Distributors(25)"
$ws.Cells.Item(9, 2).Value = 0.5441324810134711
$ws.Cells.Item(10, 1).Value = "Containers & Packaging(12)"
$ws.Cells.Item(10, 2).Value = 0.5275944728785374
$ws.Cells.Item(11, 1).Value = "Machinery(85)"
$ws.Cells.Item(11, 2).Value = 0.5214815206989146
$ws.Cells.Item(12, 1).Value = "Building Products(23)"
$ws.Cells.Item(12, 2).Value = 0.5123627667862571
$ws.Cells.Item(13, 1).Value = "Metals & Mining(89)"
$ws.Cells.Item(13, 2).Value = 0.5109642925283362
$ws.Cells.Item(14, 1).Value = "Auto Components(21)"
$ws.Cells.Item(14, 2).Value = 0.5046759631775736
$ws.Cells.Item(15, 1).Value = "Multi-Utilities(18)"
$ws.Cells.Item(15, 2).Value = 0.4909455662030818
$ws.Cells.Item(16, 1).Value = "Chemicals(51)"
$ws.Cells.Item(16, 2).Value = 0.4782659817295227
$ws.Cells.Item(17, 1).Value = "Life Sciences Tools & Services(19)"
$ws.Cells.Item(17, 2).Value = 0.4662222783206468
$ws.Cells.Item(18, 1).Value = "Wireless Telecommunication Services(14)"
$ws.Cells.Item(18, 2).Value = 0.4456355275386176
$ws.Cells.Item(19, 1).Value = "Airlines(14)"
$ws.Cells.Item(19, 2).Value = 0.4420316194123448
$ws.Cells.Item(20, 1).Value = "Insurance(75)"
$ws.Cells.Item(20, 2).Value = 0.4331316156024389
$ws.Cells.Item(21, 1).Value = "Gas Utilities(12)"
$ws.Cells.Item(21, 2).Value = 0.4197610392852889
$ws.Cells.Item(22, 1).Value = "Capital Markets(75)"
$ws.Cells.Item(22, 2).Value = 0.4179694322295037
$ws.Cells.Item(23, 1).Value = "Semiconductors & Semiconductor Equipment(68)"
$ws.Cells.Item(23, 2).Value = 0.4163578740925429
$ws.Cells.Item(24, 1).Value = "IT Services(52)"
$ws.Cells.Item(24, 2).Value = 0.4142540919728485
$ws.Cells.Item(25, 1).Value = "Leisure Products(11)"
$ws.Cells.Item(25, 2).Value = 0.4131002866899252
$ws.Cells.Item(26, 1).Value = "Electrical Equipment(28)"
$ws.Cells.Item(26, 2).Value = 0.4110709898800927
$ws.Cells.Item(27, 1).Value = "Oil, Gas & Consumable Fuels(122)"
$ws.Cells.Item(27, 2).Value = 0.4037955009463311
$ws.Cells.Item(28, 1).Value = "Household Durables(39)"
$ws.Cells.Item(28, 2).Value = 0.4036972876281053
$ws.Cells.Item(29, 1).Value = "Professional Services(35)"
$ws.Cells.Item(29, 2).Value = 0.3793475327116087
$ws.Cells.Item(30, 1).Value = "Water Utilities(12)"
$ws.Cells.Item(30, 2).Value = 0.3783444804567622
$ws.Cells.Item(31, 1).Value = "Health Care Providers & Services(46)"
$ws.Cells.Item(31, 2).Value = 0.3740882294784622
$ws.Cells.Item(32, 1).Value = "Electric Utilities(28)"
$ws.Cells.Item(32, 2).Value = 0.3626800602595638
$ws.Cells.Item(33, 1).Value = "Communications Equipment(45)"
$ws.Cells.Item(33, 2).Value = 0.3451938496260053
$ws.Cells.Item(34, 1).Value = "Banks(246)"
$ws.Cells.Item(34, 2).Value = 0.3416126992819677
$ws.Cells.Item(35, 1).Value = "Consumer Finance(15)"
$ws.Cells.Item(35, 2).Value = 0.3360127359773946
$ws.Cells.Item(36, 1).Value = "Specialty Retail(58)"
$ws.Cells.Item(36, 2).Value = 0.3341948423933782
$ws.Cells.Item(37, 1).Value = "Food & Staples Retailing(15)"
$ws.Cells.Item(37, 2).Value = 0.3328140188231028
$ws.Cells.Item(38, 1).Value = "Aerospace & Defense(37)"
$ws.Cells.Item(38, 2).Value = 0.3327326980024479
$ws.Cells.Item(39, 1).Value = "Software(66)"
$ws.Cells.Item(39, 2).Value = 0.3316468995720376
$ws.Cells.Item(40, 1).Value = "Hotels, Restaurants & Leisure(50)"
$ws.Cells.Item(40, 2).Value = 0.3270062579011052
$ws.Cells.Item(41, 1).Value = "Commercial Services & Supplies(52)"
$ws.Cells.Item(41, 2).Value = 0.3159174398305094
$ws.Cells.Item(42, 1).Value = "Textiles, Apparel & Luxury Goods(29)"
$ws.Cells.Item(42, 2).Value = 0.3084588888265488
$ws.Cells.Item(43, 1).Value = "Beverages(21)"
$ws.Cells.Item(43, 2).Value = 0.3010316192346816
$ws.Cells.Item(44, 1).Value = "Diversified Consumer Services(17)"
$ws.Cells.Item(44, 2).Value = 0.2979870840024118
$ws.Cells.Item(45, 1).Value = "Real Estate Management & Development(22)"
$ws.Cells.Item(45, 2).Value = 0.2665444198014706
$ws.Cells.Item(46, 1).Value = "Entertainment(22)"
$ws.Cells.Item(46, 2).Value = 0.2468887503039209
$ws.Cells.Item(47, 1).Value = "Media(42)"
$ws.Cells.Item(47, 2).Value = 0.2466866844775422
$ws.Cells.Item(48, 1).Value = "Diversified Telecommunication Services(20)"
$ws.Cells.Item(48, 2).Value = 0.2287517233842288
$ws.Cells.Item(49, 1).Value = "Health Care Equipment & Supplies(83)"
$ws.Cells.Item(49, 2).Value = 0.2180506714605436
$ws.Cells.Item(50, 1).Value = "Food Products(44)"
$ws.Cells.Item(50, 2).Value = 0.1978768024882402
$ws.Cells.Item(51, 1).Value = "Thrifts & Mortgage Finance(47)"
$ws.Cells.Item(51, 2).Value = 0.1785330825979096
$ws.Cells.Item(52, 1).Value = "Biotechnology(126)"
$ws.Cells.Item(52, 2).Value = 0.1760245794419086
$ws.Cells.Item(53, 1).Value = "Pharmaceuticals(48)"
$ws.Cells.Item(53, 2).Value = 0.1356876458472174
